$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.805.52"
$ws.Range("E2").Value = "  +3.02%  "
$ws.Range("D3").Value = "1.866.45"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("D4").Value = "'1.040"
$ws.Range("E4").Value = "  +3.23%  "
$ws.Range("D5").Value = "'324.93"
$ws.Range("E5").Value = "  +3.72%  "
$ws.Range("D6").Value = "'1.036"
$ws.Range("E6").Value = "  +2.89%  "
$ws.Range("D7").Value = "'0.4428"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").Value = "'0.3798"
$ws.Range("E8").Value = "  +2.87%  "
$ws.Range("D9").Value = "'0.07468"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("D10").Value = "'0.8867"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'21.79"
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("D12").Value = "1.884.56"
$ws.Range("E12").Value = "  -10.85%  "
$ws.Range("D13").Value = "'5.562"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").Value = "'6.764"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "'0.07236"
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("D16").Value = "'83.90"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "'1.041"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "'0.000009154"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").Value = "27.825.40"
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("D22").Value = "'5.340"
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("D23").Value = "'11.38"
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("D24").Value = "'1.972"
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").Value = "'158.74"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("D26").Value = "'18.90"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").Value = "'1.993"
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "'117.85"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").Value = "'0.09106"
$ws.Range("D31").Value = "'0.7813"
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("D32").Value = "'1.220"
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("D33").Value = "'3.100"
$ws.Range("E33").Value = "  +10.20%  "
$ws.Range("D34").Value = "'4.589"
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("D35").Value = "'1.038"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'0.02003"
$ws.Range("E37").Value = "  +3.91%  "
$ws.Range("D38").Value = "'0.05364"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'2.859"
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("D40").Value = "'0.5212"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").Value = "'0.1698"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").Value = "'6.896"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("D43").Value = "'8.708"
$ws.Range("E43").Value = "  +4.43%  "
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("D45").Value = "'10.73"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "'1.728"
$ws.Range("E46").Value = "  +4.55%  "
$ws.Range("D47").Value = "'0.4721"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'0.06461"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("D49").Value = "'1.903"
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("D50").Value = "'40.02"
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("D51").Value = "'64.73"
$ws.Range("E51").Value = "  +1.57%  "
